$wb = $excel.ActiveWorkbook

# --- ALERTS sheet: append two new FALL_DETECTED rows (8 and 9) ---
$alerts = $wb.Worksheets.Item("ALERTS")

$alerts.Range("A8").Value = "'2026-02-01"
$alerts.Range("B8").Value = "14:36:37"
$alerts.Range("C8").Value = "14:00"
$alerts.Range("D8").Value = "Living Room"
$alerts.Range("E8").Value = "CRITICAL"
$alerts.Range("F8").Value = "FALL_DETECTED"

$alerts.Range("A9").Value = "'2026-02-01"
$alerts.Range("B9").Value = "14:37:04"
$alerts.Range("C9").Value = "14:00"
$alerts.Range("D9").Value = "Living Room"
$alerts.Range("E9").Value = "CRITICAL"
$alerts.Range("F9").Value = "FALL_DETECTED"

# --- mmWave sheet: append one new PRESENCE_DETECTED / Active row (9) ---
$mmwave = $wb.Worksheets.Item("mmWave")

$mmwave.Range("A9").Value = "'2026-02-01"
$mmwave.Range("B9").Value = "14:36:40"
$mmwave.Range("C9").Value = "14:00"
$mmwave.Range("D9").Value = "Living Room"
$mmwave.Range("E9").Value = "PRESENCE_DETECTED"
$mmwave.Range("F9").Value = "Active"
